$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2016 -> 01/01/2023
# (shared string reused by rows 8 and 13 in both B and C columns)
# Force text format first so Excel doesn't auto-convert the date-looking
# string into a real date serial, then restore the General display format.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("B8").NumberFormat = "General"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("C8").NumberFormat = "General"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("B13").NumberFormat = "General"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("C13").NumberFormat = "General"

# Objetivos / Programa responsible person:
# 5840963 - Daniela Camargo Vernilli -> 5840897 - Clodoaldo Saron
# (shared string reused by rows 10 and 15 in both B and C columns)
$ws.Range("B10").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C10").Value = "5840897 - Clodoaldo Saron"
$ws.Range("B15").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C15").Value = "5840897 - Clodoaldo Saron"

# Critério: Avaliações escritas envolvendo -> Avaliações envolvendo
$ws.Range("B19").Value = "Avaliações envolvendo o conteúdo da disciplina."
$ws.Range("C19").Value = "Avaliações envolvendo o conteúdo da disciplina."

# Norma de recuperação: updated formula/text
$ws.Range("B20").Value = "Duas avaliações no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= média do semestre.MS> ou = 5,0: Aluno AprovadoMS< 3,0: Aluno Reprovado3,0 < ou = MS < 5,0: Aluno de Recuperação."
$ws.Range("C20").Value = "Duas avaliações no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= média do semestre.MS> ou = 5,0: Aluno AprovadoMS< 3,0: Aluno Reprovado3,0 < ou = MS < 5,0: Aluno de Recuperação."

# Bibliografia: updated recovery exam text
$ws.Range("B21").Value = "Atividade avaliativa versando sobre o conteúdo da disciplina. O aluno será aprovado se apresentar MF (média final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= média do semestre e PR= prova de recuperação."
$ws.Range("C21").Value = "Atividade avaliativa versando sobre o conteúdo da disciplina. O aluno será aprovado se apresentar MF (média final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= média do semestre e PR= prova de recuperação."
